$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Last Activity"
$ws.Range("H1").Value = "Last Welcome"
$ws.Range("G2").Value = "2025-11-12 07:46:16"
